$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.566.32'
$ws.Range('E2').Value = '  +0.11%  '

$ws.Range('D3').Value = '1.585.18'
$ws.Range('E3').Value = '  -0.47%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.80'
$ws.Range('E5').Value = '  +0.29%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.500'
$ws.Range('E6').Value = '  -0.46%  '

$ws.Range('E7').Value = '  +0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.24'
$ws.Range('E8').Value = '  -0.19%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.251'
$ws.Range('E9').Value = '  -0.52%  '

$ws.Range('E10').Value = '  -0.08%  '

$ws.Range('E11').Value = '  -0.59%  '

$ws.Range('D12').Value = '1.812.17'
$ws.Range('E12').Value = '  -0.46%  '

$ws.Range('D13').Value = '1.563.21'
$ws.Range('E13').Value = '  -1.74%  '

$ws.Range('E14').Value = '  -0.92%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.525'
$ws.Range('E15').Value = '  -2.39%  '

$ws.Range('D16').Value = '27.582.80'
$ws.Range('E16').Value = '  +0.17%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.00'
$ws.Range('E17').Value = '  -0.39%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '216.38'
$ws.Range('E18').Value = '  -1.00%  '

$ws.Range('D19').Value = '0.0₃0692'
$ws.Range('E19').Value = '  -0.42%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.31'
$ws.Range('E20').Value = '  -0.91%  '

$ws.Range('E21').Value = '  +0.15%  '

$ws.Range('E22').Value = '  -1.60%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.70'
$ws.Range('E23').Value = '  +0.18%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.00'
$ws.Range('E24').Value = '  -0.41%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.79'
$ws.Range('E25').Value = '  -1.42%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.94'
$ws.Range('E26').Value = '  +3.05%  '

$ws.Range('E27').Value = '  +0.14%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.03'
$ws.Range('E28').Value = '  +0.07%  '

$ws.Range('E29').Value = '  -1.21%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.15'
$ws.Range('E30').Value = '  -0.34%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0473'
$ws.Range('E31').Value = '  +1.12%  '

$ws.Range('E32').Value = '  -2.13%  '

$ws.Range('D33').Value = '1.372.84'
$ws.Range('E33').Value = '  +1.23%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.95'
$ws.Range('E34').Value = '  -0.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('E35').Value = '  -0.45%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.966'
$ws.Range('E36').Value = '  +0.99%  '

$ws.Range('E37').Value = '  +0.16%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0168'
$ws.Range('E38').Value = '  +1.40%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.533'
$ws.Range('E39').Value = '  -1.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.822'
$ws.Range('E40').Value = '  +1.21%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.972'
$ws.Range('E42').Value = '  +0.38%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '64.15'
$ws.Range('E43').Value = '  +0.36%  '

$ws.Range('E44').Value = '  +4.19%  '

$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.26'
$ws.Range('E45').Value = '  -1.59%  '

$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.76'
$ws.Range('E46').Value = '  -0.72%  '

$ws.Range('D47').Value = '1.722.83'
$ws.Range('E47').Value = '  -0.40%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.40'
$ws.Range('E48').Value = '  -1.66%  '

$ws.Range('E49').Value = '  +6.00%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0958'
$ws.Range('E50').Value = '  -1.22%  '

$ws.Range('E51').Value = '  -0.72%  '
